$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the "Run By" row (row 10). This pushes
# "Run By" -> row 11, "Run Date" -> row 12, the blank spacer row -> row 13,
# and the column-header row -> row 14 (matching the new dimension A1:BB14).
$ws.Rows.Item(10).Insert()

# The new row 10 should look like the "Date From / Date To" row above it
# (two label+input pairs sharing the same formatting), so copy that
# formatting down before typing the new "Charge Type" filter label.
$ws.Range("A9:D9").Copy()
$ws.Range("A10:D10").PasteSpecial(-4122)

# Only column A carries the visible "Charge Type" label; B/C/D stay blank
# placeholders (matching the template's unused second label/input slot).
$ws.Range("A10").Value = "Charge Type"
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""

$excel.CutCopyMode = 0
